$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill H2:H6 with the password value used by userRegistration first so it
# lands earlier in the shared-strings table, then add the "pass" header.
$ws.Range("H2:H6").Value = "blackdress19"
$ws.Range("H1").Value = "pass"

# Adjust column widths to match target layout (internal width = ColumnWidth + 5/6)
$ws.Columns.Item(7).ColumnWidth = 5.166666666666667
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666

# Update the active selection to reflect where the user left off editing
$ws.Range("H13").Select()
